# SCSG Quarterly Financials update
# Insert two new quarter columns (newest quarters) in front of the existing
# "Period Ending" column D, shifting all the historical quarters right, and
# populate the new columns with the latest quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert 2 new columns before column D. This shifts old D:K -> F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) Re-apply number formatting to the freshly inserted D:E columns so they
#    match the columns that used to live there (date format for the
#    "Period Ending" header rows, thousands-format for the data rows).
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $ws.Range("F$r`:G$r").Copy()
    $ws.Range("D$r`:E$r").PasteSpecial(-4122)
}

$numberBlocks = @(@(8,35), @(39,77), @(81,102))
foreach ($block in $numberBlocks) {
    $startRow = $block[0]
    $endRow = $block[1]
    $ws.Range("F$startRow`:G$endRow").Copy()
    $ws.Range("D$startRow`:E$endRow").PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = 0

# 3) Fill in the new column D (latest quarter) / E (previous quarter) values.
$ws.Range("D7").Value = 43373; $ws.Range("E7").Value = 43281
$ws.Range("D8").Value = 5300; $ws.Range("E8").Value = 5100
$ws.Range("D9").Value = "NA"; $ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"; $ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"; $ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0; $ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0; $ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0; $ws.Range("E15").Value = 0
$ws.Range("D17").Value = 600; $ws.Range("E17").Value = 500
$ws.Range("D18").Value = 4700; $ws.Range("E18").Value = 4600
$ws.Range("D20").Value = -3200; $ws.Range("E20").Value = -3600
$ws.Range("D21").Value = "NA"; $ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0; $ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1400; $ws.Range("E23").Value = 1000
$ws.Range("D24").Value = 300; $ws.Range("E24").Value = 200
$ws.Range("D25").Value = 0; $ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1100; $ws.Range("E26").Value = 800
$ws.Range("D27").Value = 1100; $ws.Range("E27").Value = 800
$ws.Range("D28").Value = 0; $ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"; $ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0; $ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0; $ws.Range("E31").Value = 0
$ws.Range("D32").Value = 3200; $ws.Range("E32").Value = 3600
$ws.Range("D33").Value = 1100; $ws.Range("E33").Value = 800
$ws.Range("D34").Value = 0; $ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1100; $ws.Range("E35").Value = 800
$ws.Range("D38").Value = 43373; $ws.Range("E38").Value = 43281
$ws.Range("D41").Value = 14200; $ws.Range("E41").Value = 18400
$ws.Range("D42").Value = 1300; $ws.Range("E42").Value = 6800
$ws.Range("D43").Value = 0; $ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0; $ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0; $ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0; $ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0; $ws.Range("E47").Value = 0
$ws.Range("D48").Value = 9700; $ws.Range("E48").Value = 10000
$ws.Range("D49").Value = 100; $ws.Range("E49").Value = 200
$ws.Range("D50").Value = 0; $ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0; $ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0; $ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0; $ws.Range("E53").Value = 0
$ws.Range("D54").Value = 544100; $ws.Range("E54").Value = 541300
$ws.Range("D57").Value = 0; $ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0; $ws.Range("E58").Value = 0
$ws.Range("D59").Value = "NA"; $ws.Range("E59").Value = "NA"
$ws.Range("D60").Value = 0; $ws.Range("E60").Value = 0
$ws.Range("D61").Value = 50600; $ws.Range("E61").Value = 40000
$ws.Range("D62").Value = 0; $ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0; $ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0; $ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0; $ws.Range("E65").Value = 0
$ws.Range("D66").Value = 491300; $ws.Range("E66").Value = 488400
$ws.Range("D68").Value = 0; $ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0; $ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0; $ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0; $ws.Range("E71").Value = 0
$ws.Range("D72").Value = "NA"; $ws.Range("E72").Value = "NA"
$ws.Range("D73").Value = 0; $ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0; $ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0; $ws.Range("E75").Value = 0
$ws.Range("D76").Value = 52800; $ws.Range("E76").Value = 52900
$ws.Range("D77").Value = 0; $ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43373; $ws.Range("E80").Value = 43281
$ws.Range("D81").Value = 1100; $ws.Range("E81").Value = 800
$ws.Range("D83").Value = 0; $ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0; $ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0; $ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0; $ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0; $ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0; $ws.Range("E88").Value = 0
$ws.Range("D89").Value = 0; $ws.Range("E89").Value = 0
$ws.Range("D91").Value = 0; $ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0; $ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0; $ws.Range("E93").Value = 0
$ws.Range("D94").Value = 0; $ws.Range("E94").Value = 0
$ws.Range("D96").Value = 0; $ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0; $ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0; $ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0; $ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0; $ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0; $ws.Range("E101").Value = 0
$ws.Range("D102").Value = 0; $ws.Range("E102").Value = 0
